# Insert a new weekly price record for "Damasco" (Modesto variety) at row 56
# of the "Vega Modelo de Temuco" sheet, shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (shifts rows 56:80 down to 57:81)
$ws.Rows(56).Insert()

# Fill in the new row with the new record's data
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44917
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = "Damasco"
$ws.Range("K56").Value = "Modesto"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 400
$ws.Range("N56").Value = 20000
$ws.Range("O56").Value = 22000
$ws.Range("P56").Value = 20500
$ws.Range("Q56").Value = "`$/bandeja 18 kilos"
$ws.Range("R56").Value = "Región de O'Higgins"
$ws.Range("S56").Value = 1139
$ws.Range("T56").Value = 18
